$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1246.0714
$ws.Range("I6").Value = 1265
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 3795
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -3683
$ws.Range("N6").Value = -3224

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3685
$ws.Range("N76").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2908
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 104174080
$ws.Range("I86").Value = 111118680
$ws.Range("K86").Value = 111118680
$ws.Range("M86").Value = -111117557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 104174080
$ws.Range("I89").Value = 111118680
$ws.Range("K89").Value = 555593400
$ws.Range("M89").Value = -555587784

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 181821740
$ws.Range("I106").Value = 181821740
$ws.Range("K106").Value = 181821740
$ws.Range("M106").Value = -181821109

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1161.1666
$ws.Range("I135").Value = 1161.1666
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 10450.4994
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -7915.499400000001
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1305013.9
$ws.Range("I137").Value = 6375.5454
$ws.Range("J137").Value = 2278992.8
$ws.Range("K137").Value = 19126.6362
$ws.Range("L137").Value = 6836978.399999999
$ws.Range("M137").Value = -16576.6362
$ws.Range("N137").Value = -6842078.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 38000
$ws.Range("I15").Value = 27500
$ws.Range("J15").Value = 45000
$ws.Range("K15").Value = 27500
$ws.Range("L15").Value = 45000
$ws.Range("M15").Value = -27150
$ws.Range("N15").Value = -45700

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3763.64
$ws.Range("I32").Value = 2204
$ws.Range("K32").Value = 2204
$ws.Range("M32").Value = -1917

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1876036.6
$ws.Range("I61").Value = 59086.74
$ws.Range("J61").Value = 4531579
$ws.Range("K61").Value = 59086.74
$ws.Range("L61").Value = 4531579
$ws.Range("M61").Value = -58874.74
$ws.Range("N61").Value = -4532003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 149571.28
$ws.Range("J86").Value = 149571.28
$ws.Range("L86").Value = 149571.28
$ws.Range("N86").Value = -151943.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 149571.28
$ws.Range("J89").Value = 149571.28
$ws.Range("L89").Value = 448713.84
$ws.Range("N89").Value = -460569.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1876036.6
$ws.Range("I136").Value = 59086.74
$ws.Range("J136").Value = 4531579
$ws.Range("K136").Value = 177260.22
$ws.Range("L136").Value = 13594737
$ws.Range("M136").Value = -174710.22
$ws.Range("N136").Value = -13599837

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 31125.2
$ws.Range("J103").Value = 31125.2
$ws.Range("L103").Value = 31125.2
$ws.Range("N103").Value = -33469.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 36069
$ws.Range("J114").Value = 36069
$ws.Range("L114").Value = 36069
$ws.Range("N114").Value = -44747

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 100000
$ws.Range("J117").Value = 100000
$ws.Range("L117").Value = 100000
$ws.Range("N117").Value = -109178

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 30002412
$ws.Range("I134").Value = 2406.0715
$ws.Range("K134").Value = 7218.2145
$ws.Range("M134").Value = -4683.2145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 6499.5
$ws.Range("I19").Value = 499.5
$ws.Range("J19").Value = 12499.5
$ws.Range("K19").Value = 499.5
$ws.Range("L19").Value = 12499.5
$ws.Range("M19").Value = -329.5
$ws.Range("N19").Value = -12839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 6499.5
$ws.Range("I24").Value = 499.5
$ws.Range("J24").Value = 12499.5
$ws.Range("K24").Value = 499.5
$ws.Range("L24").Value = 12499.5
$ws.Range("M24").Value = -329.5
$ws.Range("N24").Value = -12839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4535.3853
$ws.Range("I31").Value = 3107
$ws.Range("J31").Value = 4842.7593
$ws.Range("K31").Value = 3107
$ws.Range("L31").Value = 4842.7593
$ws.Range("M31").Value = -2812
$ws.Range("N31").Value = -5432.7593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4535.3853
$ws.Range("I34").Value = 3107
$ws.Range("J34").Value = 4842.7593
$ws.Range("K34").Value = 3107
$ws.Range("L34").Value = 4842.7593
$ws.Range("M34").Value = -2905
$ws.Range("N34").Value = -5246.7593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2566.6
$ws.Range("I58").Value = 2250
$ws.Range("J58").Value = 2928.4285
$ws.Range("K58").Value = 2250
$ws.Range("L58").Value = 2928.4285
$ws.Range("M58").Value = -2047
$ws.Range("N58").Value = -3334.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 66371
$ws.Range("J116").Value = 66371
$ws.Range("L116").Value = 66371
$ws.Range("N116").Value = -75549

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2566.6
$ws.Range("I136").Value = 2250
$ws.Range("J136").Value = 2928.4285
$ws.Range("K136").Value = 6750
$ws.Range("L136").Value = 8785.2855
$ws.Range("M136").Value = -4200
$ws.Range("N136").Value = -13885.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 4666.5
$ws.Range("J20").Value = 4750
$ws.Range("L20").Value = 14250
$ws.Range("N20").Value = -14704

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 20995
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 20995
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 62985
$ws.Range("N82").Value = -63797
$ws.Range("M82").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 20995
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 20995
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 62985
$ws.Range("N85").Value = -65793
$ws.Range("M85").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 513430.47
$ws.Range("J94").Value = 627057.3
$ws.Range("L94").Value = 1881171.9
$ws.Range("N94").Value = -1882523.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 912124.25
$ws.Range("I103").Value = 1429195.4
$ws.Range("J103").Value = 7249.75
$ws.Range("K103").Value = 4287586.199999999
$ws.Range("L103").Value = 21749.25
$ws.Range("M103").Value = -4286707.199999999
$ws.Range("N103").Value = -23507.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 8425.6
$ws.Range("J106").Value = 8425.6
$ws.Range("L106").Value = 25276.8
$ws.Range("N106").Value = -27168.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 83349550
$ws.Range("I117").Value = 3233.3333
$ws.Range("J117").Value = 133357336
$ws.Range("K117").Value = 9699.999899999999
$ws.Range("L117").Value = 400072008
$ws.Range("M117").Value = -6257.999899999999
$ws.Range("N117").Value = -400078892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 11807.4
$ws.Range("I120").Value = 9764.25
$ws.Range("J120").Value = 19980
$ws.Range("K120").Value = 29292.75
$ws.Range("L120").Value = 59940
$ws.Range("M120").Value = -24454.75
$ws.Range("N120").Value = -69616

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 3869.3057
$ws.Range("I124").Value = 2431.6667
$ws.Range("J124").Value = 4000
$ws.Range("K124").Value = 7295.000100000001
$ws.Range("L124").Value = 12000
$ws.Range("M124").Value = -2385.000100000001
$ws.Range("N124").Value = -21820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 7487.3477
$ws.Range("J129").Value = 11357.643
$ws.Range("L129").Value = 34072.929
$ws.Range("N129").Value = -44072.929

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3954162.8
$ws.Range("I131").Value = 7577036.5
$ws.Range("J131").Value = 1937
$ws.Range("K131").Value = 22731109.5
$ws.Range("L131").Value = 5811
$ws.Range("M131").Value = -22726069.5
$ws.Range("N131").Value = -15891

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 3499
$ws.Range("I133").Value = 3499
$ws.Range("K133").Value = 10497
$ws.Range("M133").Value = -5437

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2320.75
$ws.Range("I137").Value = 1616.6666
$ws.Range("J137").Value = 4433
$ws.Range("K137").Value = 4849.9998
$ws.Range("L137").Value = 13299
$ws.Range("M137").Value = 250.0002000000004
$ws.Range("N137").Value = -23499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 14999.5
$ws.Range("I20").Value = 14999.5
$ws.Range("K20").Value = 14999.5
$ws.Range("M20").Value = -14754.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 60000
$ws.Range("J114").Value = 60000
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 16981.334
$ws.Range("I4").Value = 16378.2
$ws.Range("K4").Value = 16378.2
$ws.Range("M4").Value = -16265.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 16981.334
$ws.Range("I28").Value = 16378.2
$ws.Range("K28").Value = 16378.2
$ws.Range("M28").Value = -16146.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 16981.334
$ws.Range("I37").Value = 16378.2
$ws.Range("K37").Value = 16378.2
$ws.Range("M37").Value = -16271.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 54422.4
$ws.Range("I136").Value = 63008.766
$ws.Range("J136").Value = 5766.3335
$ws.Range("K136").Value = 189026.298
$ws.Range("L136").Value = 17299.0005
$ws.Range("M136").Value = -186476.298
$ws.Range("N136").Value = -22399.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3584.158
$ws.Range("I136").Value = 3153.25
$ws.Range("J136").Value = 4322.857
$ws.Range("K136").Value = 9459.75
$ws.Range("L136").Value = 12968.571
$ws.Range("M136").Value = -6909.75
$ws.Range("N136").Value = -18068.571
